$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: 特殊物品编号 / 特殊物品名称 / 获取途径
$ws.Range("A1").Value = "特殊物品编号"
$ws.Range("B1").Value = "特殊物品名称"
$ws.Range("C1").Value = "获取途径"

# Column widths (as authored in the source workbook)
$ws.Columns.Item(1).ColumnWidth = 16.3
$ws.Columns.Item(2).ColumnWidth = 15.15
$ws.Columns.Item(3).ColumnWidth = 14.3

# Restore the selection left behind by the editing session
$ws.Range("G5").Select() | Out-Null
